$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column J (DIAS) for rows 2 through 29: each value increases by 4
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 10)  # column J = 10
    $cell.Value2 = $cell.Value2 + 4
}
